$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.430.44'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '2.585.07'
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '591.27'
$ws.Range("E5").Value = '  -2.72%  '
$ws.Range("D6").Value = '150.92'
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("D10").Value = '5.71'
$ws.Range("E10").Value = '  +2.04%  '
$ws.Range("D11").Value = '0.385'
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("D13").Value = '27.61'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").Value = '3.046.75'
$ws.Range("D15").Value = '63.195.43'
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("E16").Value = '  +5.82%  '
$ws.Range("D17").Value = '2.575.82'
$ws.Range("E17").Value = '  -2.96%  '
$ws.Range("D18").Value = '12.27'
$ws.Range("E18").Value = '  +4.08%  '
$ws.Range("D19").Value = '4.74'
$ws.Range("E19").Value = '  +3.93%  '
$ws.Range("D20").Value = '345.55'
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = '67.45'
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").Value = '1.69'
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("D25").Value = '9.25'
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").Value = '563.74'
$ws.Range("E27").Value = '  -0.47%  '
$ws.Range("D28").Value = '8.08'
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = '2.04'
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").Value = '0.0₃0852'
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("D34").Value = '5.25'
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("D35").Value = '166.51'
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("D36").Value = '0.413'
$ws.Range("E36").Value = '  +2.01%  '
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = '19.55'
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  -1.33%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '166.25'
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("D42").Value = '39.57'
$ws.Range("E42").Value = '  -1.46%  '
$ws.Range("D43").Value = '3.99'
$ws.Range("E43").Value = '  +5.32%  '
$ws.Range("D44").Value = '22.90'
$ws.Range("E44").Value = '  +4.46%  '
$ws.Range("D45").Value = '0.0585'
$ws.Range("E45").Value = '  +3.38%  '
$ws.Range("D46").Value = '2.11'
$ws.Range("E46").Value = '  +5.42%  '
$ws.Range("D47").Value = '0.629'
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("E48").Value = '  +2.81%  '
$ws.Range("D49").Value = '0.0962'
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("D50").Value = '19.15'
$ws.Range("E50").Value = '  +1.61%  '
$ws.Range("D51").Value = '0.0₆0232'
$ws.Range("E51").Value = '  +17.61%  '
